$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string / label text updates ---
# A8 is a rich-text run "Volume " / "30" / "   Number  " / "27" -> only the
# trailing issue number run changes from "27" to "28".
$volRange = $ws.Range("A8")
$volRange.Characters(21, 2).Text = "28"

# C9 is a rich-text run "Report Covering the Week  " / "7/3/2023" /
# "  Through  " / "7/9/2023" -> both date runs change.
$weekRange = $ws.Range("C9")
$weekRange.Characters(27, 8).Text = "7/10/2023"
$weekRange.Characters(47, 8).Text = "7/16/2023"

# --- Crime Complaints table numeric updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = -80
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = -21.428571428571
$ws.Range("I14").Value = 71
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = -11.25
$ws.Range("L14").Value = -5.333333333333
$ws.Range("M14").Value = 4.411764705882
$ws.Range("N14").Value = -73.80073800738
# Row 15
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 28
$ws.Range("G15").Value = 32
$ws.Range("H15").Value = -12.5
$ws.Range("I15").Value = 212
$ws.Range("J15").Value = 218
$ws.Range("K15").Value = -2.752293577981
$ws.Range("L15").Value = 12.765957446808
$ws.Range("M15").Value = 39.473684210526
$ws.Range("N15").Value = -44.356955380577
# Row 16
$ws.Range("D16").Value = 146
$ws.Range("E16").Value = -30.821917808219
$ws.Range("F16").Value = 440
$ws.Range("G16").Value = 503
$ws.Range("H16").Value = -12.524850894632
$ws.Range("I16").Value = 2490
$ws.Range("J16").Value = 2651
$ws.Range("K16").Value = -6.073179932101
$ws.Range("L16").Value = 33.583690987124
$ws.Range("M16").Value = 9.450549450549
$ws.Range("N16").Value = -71.039776692254
# Row 17
$ws.Range("C17").Value = 163
$ws.Range("D17").Value = 166
$ws.Range("E17").Value = -1.807228915662
$ws.Range("F17").Value = 725
$ws.Range("G17").Value = 662
$ws.Range("H17").Value = 9.516616314199
$ws.Range("I17").Value = 4226
$ws.Range("J17").Value = 3921
$ws.Range("K17").Value = 7.778627901045
$ws.Range("L17").Value = 32.269170579029
$ws.Range("M17").Value = 79.60050998725
$ws.Range("N17").Value = -13.913220615196
# Row 18
$ws.Range("C18").Value = 63
$ws.Range("D18").Value = 52
$ws.Range("E18").Value = 21.153846153846
$ws.Range("F18").Value = 210
$ws.Range("G18").Value = 238
$ws.Range("H18").Value = -11.764705882352
$ws.Range("I18").Value = 1607
$ws.Range("J18").Value = 1604
$ws.Range("K18").Value = 0.187032418952
$ws.Range("L18").Value = 44.254937163375
$ws.Range("M18").Value = -6.187974314068
$ws.Range("N18").Value = -84.017901541521
# Row 19
$ws.Range("C19").Value = 152
$ws.Range("D19").Value = 163
$ws.Range("E19").Value = -6.748466257668
$ws.Range("F19").Value = 593
$ws.Range("G19").Value = 636
$ws.Range("H19").Value = -6.761006289308
$ws.Range("I19").Value = 4116
$ws.Range("J19").Value = 4246
$ws.Range("K19").Value = -3.061705134243
$ws.Range("L19").Value = 23.78947368421
$ws.Range("M19").Value = 71.71464330413
$ws.Range("N19").Value = 5.187835420393
# Row 20
$ws.Range("C20").Value = 102
$ws.Range("D20").Value = 86
$ws.Range("E20").Value = 18.60465116279
$ws.Range("F20").Value = 376
$ws.Range("G20").Value = 323
$ws.Range("H20").Value = 16.40866873065
$ws.Range("I20").Value = 2851
$ws.Range("J20").Value = 2165
$ws.Range("K20").Value = 31.685912240184
$ws.Range("L20").Value = 106.744017403916
$ws.Range("M20").Value = 162.523020257827
$ws.Range("N20").Value = -65.303638797614
# Row 21
$ws.Range("C21").Value = 590
$ws.Range("D21").Value = 626
$ws.Range("E21").Value = -5.750798722044
$ws.Range("F21").Value = 2383
$ws.Range("G21").Value = 2408
$ws.Range("H21").Value = -1.038205980066
$ws.Range("I21").Value = 15573
$ws.Range("J21").Value = 14885
$ws.Range("K21").Value = 4.622102788041
$ws.Range("L21").Value = 39.793536804308
$ws.Range("M21").Value = 55.047789725209
$ws.Range("N21").Value = -57.151111600264
# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -28.571428571428
$ws.Range("F22").Value = 16
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = -46.666666666666
$ws.Range("I22").Value = 158
$ws.Range("J22").Value = 198
$ws.Range("K22").Value = -20.20202020202
$ws.Range("L22").Value = 27.419354838709
$ws.Range("M22").Value = -9.19540229885
# Row 23
$ws.Range("C23").Value = 36
$ws.Range("D23").Value = 34
$ws.Range("E23").Value = 5.882352941176
$ws.Range("F23").Value = 134
$ws.Range("H23").Value = -6.293706293706
$ws.Range("I23").Value = 949
$ws.Range("J23").Value = 877
$ws.Range("K23").Value = 8.209806157354
$ws.Range("L23").Value = 49.68454258675
$ws.Range("M23").Value = 68.262411347517
# Row 24
$ws.Range("C24").Value = 354
$ws.Range("D24").Value = 381
$ws.Range("E24").Value = -7.086614173228
$ws.Range("F24").Value = 1421
$ws.Range("G24").Value = 1623
$ws.Range("H24").Value = -12.446087492298
$ws.Range("I24").Value = 9583
$ws.Range("J24").Value = 9893
$ws.Range("K24").Value = -3.133528757707
$ws.Range("L24").Value = 46.216051266402
$ws.Range("M24").Value = 43.565543071161
# Row 25
$ws.Range("C25").Value = 209
$ws.Range("D25").Value = 184
$ws.Range("E25").Value = 13.586956521739
$ws.Range("F25").Value = 852
$ws.Range("G25").Value = 892
$ws.Range("H25").Value = -4.484304932735
$ws.Range("I25").Value = 5689
$ws.Range("J25").Value = 5467
$ws.Range("K25").Value = 4.060728004389
$ws.Range("L25").Value = 28.448859787762
$ws.Range("M25").Value = -5.10425354462
# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -30.76923076923
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 59
$ws.Range("H26").Value = -33.898305084745
$ws.Range("I26").Value = 346
$ws.Range("J26").Value = 384
$ws.Range("K26").Value = -9.895833333333
$ws.Range("L26").Value = 11.254019292604
# Row 27
$ws.Range("C27").Value = 21
$ws.Range("D27").Value = 21
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 80
$ws.Range("G27").Value = 71
$ws.Range("H27").Value = 12.676056338028
$ws.Range("I27").Value = 570
$ws.Range("J27").Value = 491
$ws.Range("K27").Value = 16.089613034623
$ws.Range("L27").Value = 21.535181236673
# Row 28
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 15
$ws.Range("E28").Value = -40
$ws.Range("F28").Value = 62
$ws.Range("G28").Value = 60
$ws.Range("H28").Value = 3.333333333333
$ws.Range("I28").Value = 219
$ws.Range("J28").Value = 293
$ws.Range("K28").Value = -25.255972696245
$ws.Range("L28").Value = -32.82208588957
$ws.Range("M28").Value = -10.245901639344
$ws.Range("N28").Value = -69.456066945606
# Row 29
$ws.Range("C29").Value = 6
$ws.Range("E29").Value = -53.846153846153
$ws.Range("F29").Value = 47
$ws.Range("G29").Value = 51
$ws.Range("H29").Value = -7.843137254901
$ws.Range("I29").Value = 179
$ws.Range("J29").Value = 251
$ws.Range("K29").Value = -28.685258964143
$ws.Range("L29").Value = -35.144927536231
$ws.Range("M29").Value = -11.822660098522
$ws.Range("N29").Value = -72.376543209876
# Row 30
$ws.Range("J30").Value = 28
$ws.Range("K30").Value = -57.142857142857
